# Applies the "Week 02" sheet build-out + related workbook/view tweaks
# described by the authoritative diff.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Week 01")
$ws3 = $wb.Worksheets.Item("Week 02")

# --- 1. Clone the "Week 01" layout (values + merges, then formats) into
#        "Week 02" so column widths, borders, fonts, and merged ranges all
#        match the existing weekly template. ---
$ws2.Range("A1:R9").Copy()
$ws3.Range("A1").PasteSpecial()
$ws2.Range("A1:R9").Copy()
$ws3.Range("A1").PasteSpecial(-4122)

# --- 2. Insert the new member row (row 7, "Nguyen Quang Vinh") pushing the
#        previous row 7 ("Thang") and the blank rows down by one, then copy
#        the sibling row's formatting onto the freshly inserted row. ---
$ws3.Rows.Item(7).Insert()
$ws3.Range("B8:R8").Copy()
$ws3.Range("B7").PasteSpecial(-4122)

# --- 3. Row heights (row 5 grows to fit the long wrapped comment; the rest
#        keep the template's 15pt / 5.25pt heights that PasteSpecial already
#        carried along implicitly via row 8/9 data, but re-assert them so
#        the freshly inserted / copied rows line up exactly). ---
$ws3.Rows.Item(1).RowHeight = 5.25
$ws3.Rows.Item(2).RowHeight = 20.4
$ws3.Rows.Item(3).RowHeight = 16.2
$ws3.Rows.Item(4).RowHeight = 15
$ws3.Rows.Item(5).RowHeight = 40.8
$ws3.Rows.Item(6).RowHeight = 15
$ws3.Rows.Item(7).RowHeight = 15
$ws3.Rows.Item(8).RowHeight = 15
$ws3.Rows.Item(9).RowHeight = 15
$ws3.Rows.Item(10).RowHeight = 15

# --- 4. Title / header text for the new week. ---
$ws3.Range("B2").Value = "9h00-9h15 ngày 17/9/2024"

# --- 5. New member row content (row 7). ---
$ws3.Range("B7").Value = "Nguyễn Quang Vinh"
$ws3.Range("C7").Value = 1
$ws3.Range("I7").Value = 1
$ws3.Range("L7").ClearContents()
$ws3.Range("M7").Value = 1
$ws3.Range("R7").Value = 2

# --- 6. "Da lam" (done), "Kho khan" (difficulties) and "Nhan xet" (remarks)
#        columns for every member row. ---
$ws3.Range("F5").Value = "Nguyễn Minh Quân quay buổi họp, phân việc "
$ws3.Range("G5").Value = "đang hơi chậm tiến độ"
$ws3.Range("Q5").Value = "cân quy định giờ giấc hợp lí và tổ chức buổi họp thường xuyên"

$ws3.Range("F6").Value = "Chung vẽ giao diện FIGMA"
$ws3.Range("G6").Value = "đang hơi chậm tiến độ"
$ws3.Range("Q6").Value = "xây dựng có tích cực "

$ws3.Range("F7").Value = "Vinh làm sơ đồ ERF"
$ws3.Range("G7").Value = "đang hơi chậm tiến độ"
$ws3.Range("Q7").Value = "nên xây dựng hơn chút "

$ws3.Range("F8").Value = "Thắng làm User Story và Sprint Backlog "
$ws3.Range("G8").Value = "đang hơi chậm tiến độ"
$ws3.Range("Q8").Value = "tích cực"

# --- 7. View bits: "Week 02" becomes the active / selected tab, scrolled to
#        column G, with H13 selected; "Week 01" loses its former
#        tab-selected / topLeftCell state and ends up with a whole-sheet
#        selection (its last interactive state before focus moved on). ---
$ws2.Activate()
$ws2.Cells.Select()

$ws3.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws3.Range("H13").Select()
